$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All affected cells are plain text (inlineStr) cells in the source workbook.
# Force text format first so Excel does not reinterpret numeric-looking strings
# (e.g. "47.349.37", "0.140") as numbers/dates when the value is assigned.
$cellRefs = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "E9", "D10", "E10", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "E36", "D37", "E37", "E38", "E39", "B40", "C40", "D40", "E40", "B41", "C41", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "B48", "C48", "D48", "E48", "B49", "C49", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '47.349.37'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '2.488.96'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '321.54'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = '108.90'
$ws.Range("E6").Value = '  +3.32%  '
$ws.Range("D7").Value = '0.521'
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").Value = '39.41'
$ws.Range("E10").Value = '  +4.32%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '18.57'
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("D14").Value = '7.19'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '2.877.67'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '2.491.01'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '0.846'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").Value = '47.255.05'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("E19").Value = '  +6.07%  '
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("D21").Value = '0.0₃0940'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("E22").Value = '  +15.50%  '
$ws.Range("D23").Value = '70.61'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '246.94'
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").Value = '2.54'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '25.71'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("E28").Value = '  +3.93%  '
$ws.Range("D29").Value = '9.96'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  +4.40%  '
$ws.Range("D31").Value = '34.52'
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").Value = '49.84'
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").Value = '20.39'
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("D34").Value = '5.30'
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("D35").Value = '0.0786'
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '4.75'
$ws.Range("E37").Value = '  +2.65%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.112'
$ws.Range("E40").Value = '  +0.42%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '22.53'
$ws.Range("E41").Value = '  +7.82%  '
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = '119.06'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").Value = '0.0296'
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").Value = '1.992.17'
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("E46").Value = '  +2.02%  '
$ws.Range("E47").Value = '  -3.02%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '1.78'
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '9.07'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").Value = '5.19'
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("D51").Value = '56.77'
$ws.Range("E51").Value = '  +3.65%  '
